$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update task descriptions:
# - Massine Merouane (row 5) now has the task "Front"
# - Benoit Danglades (row 3) and Kevin Amadji (row 4) now share the task "BDD+Back"
$ws.Range("C5").Value = "Front"
$ws.Range("C3").Value = "BDD+Back"
$ws.Range("C4").Value = "BDD+Back"

# Update the active selection to C6
$ws.Range("C6").Select()
